# This script re-orders the data rows (rows 2-32) of the active worksheet.
# Each target row ends up containing the values that used to live in a
# different row (a single 31-element permutation cycle). Row 1 (header)
# and row 33 are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 32
$firstCol = 1   # A
$lastCol = 51   # AY

# Mapping: target row -> source row (i.e. after the edit, row <key> holds
# the data that, before the edit, lived in row <value>).
$rowMap = @{
    2  = 15
    3  = 16
    4  = 17
    5  = 2
    6  = 3
    7  = 18
    8  = 4
    9  = 19
    10 = 5
    11 = 6
    12 = 20
    13 = 21
    14 = 22
    15 = 7
    16 = 8
    17 = 9
    18 = 23
    19 = 24
    20 = 25
    21 = 26
    22 = 10
    23 = 27
    24 = 28
    25 = 11
    26 = 12
    27 = 29
    28 = 30
    29 = 13
    30 = 31
    31 = 32
    32 = 14
}

# Columns that hold text which Excel's automatic type detection would
# otherwise reinterpret as a number or a date when assigned through
# Value2 (e.g. "1" -> 1, "2023-08-28" -> a date serial). Force those
# columns to a plain text number format before writing so the original
# text representation is preserved.
$textColumns = @(9, 25, 27)   # I, Y, AA

# 1. Snapshot the full block of data (rows 2-32, columns A-AY) before
#    making any changes, so overwriting rows doesn't clobber data that
#    still needs to be read for a later row in the permutation.
$srcRange = $ws.Range($ws.Cells.Item($firstRow, $firstCol), $ws.Cells.Item($lastRow, $lastCol))
$snapshot = $srcRange.Value2

# 2. Make sure the text columns keep their text representation.
foreach ($col in $textColumns) {
    $colRange = $ws.Range($ws.Cells.Item($firstRow, $col), $ws.Cells.Item($lastRow, $col))
    $colRange.NumberFormat = "@"
}

# 3. Write each target row using the snapshot data taken from its source
#    row.
foreach ($targetRow in ($rowMap.Keys | Sort-Object)) {
    $sourceRow = $rowMap[$targetRow]
    $snapshotRowIndex = $sourceRow - $firstRow + 1

    $rowValues = New-Object 'object[,]' 1, ($lastCol - $firstCol + 1)
    for ($c = 1; $c -le ($lastCol - $firstCol + 1); $c++) {
        $rowValues[0, $c - 1] = $snapshot[$snapshotRowIndex, $c]
    }

    $destRange = $ws.Range($ws.Cells.Item($targetRow, $firstCol), $ws.Cells.Item($targetRow, $lastCol))
    $destRange.Value2 = $rowValues
}
